# Apply "Updated group list and fixed spelling mistakes" edits.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Group List")

# --- Fix spelling mistakes in header row (row 1) ---
$ws.Range("G1").Value = "Strength"
$ws.Range("H1").Value = "Weakness"
$ws.Range("I1").Value = "Lead"

# --- Row 2 (Ben Dudley): update Strength / Lead columns ---
$ws.Range("G2").Value = "DB/Basic Java"
$ws.Range("I2").Value = "DB (y)"

# --- Row 7 (Maurice Corriette): add blog links + fill in Strength/Weakness/Lead ---
$ws.Hyperlinks.Add($ws.Range("E7"), "http://mac81cs.blogspot.co.uk/")
$ws.Range("E7").Style = $ws.Range("E2").Style
$ws.Range("F7").Value = "http://mac81cs.blogspot.com/feeds/posts/default"
$ws.Range("G7").Value = "Testing/Web/DB"
$ws.Range("H7").Value = "Java"
$ws.Range("I7").Value = "-"

# --- Row 8 (Oliver Earl): fill in Lead column ---
$ws.Range("I8").Value = "Web"

# --- Row 9 (Tim Anderson): update Strength / Lead columns ---
$ws.Range("G9").Value = "Web/Basic Java"
$ws.Range("I9").Value = "-"

# --- Update selection / scroll position ---
$ws.Range("A4").Select()
